$wb = $excel.ActiveWorkbook

$updates = @{
    "Vendor Not Found"      = @{ "A2" = "Vendor_1663510155044" }
    "Microsoft"             = @{ "A2" = "Microsoft_1663510154692"; "A3" = "Microsoft_1663510154970"; "A4" = "Microsoft_1663510155345" }
    "Citrix"                = @{ "A2" = "Citrix_1663510154851"; "A3" = "Citrix_1663510154901"; "A4" = "Citrix_1663510154935"; "A5" = "Citrix_1663510155312"; "A6" = "Citrix_1663510155486" }
    "ConnectWise"           = @{ "A2" = "ConnectWise_1663510155589" }
    "CISCO"                 = @{ "A2" = "CISCO_1663510154750"; "A3" = "CISCO_1663510155618" }
    "SAP"                   = @{ "A2" = "SAP_1663510155645" }
    "Fortinet Technologies" = @{ "A2" = "Fortinet_1663510154817" }
    "Oracle"                = @{ "A2" = "Oracle_1663510155011"; "A3" = "Oracle_1663510155384" }
    "Test Vendor 1"         = @{ "A2" = "Test_1663510155080" }
    "Test Vendor 2"         = @{ "A2" = "Test_1663510155127" }
    "Waves Audio"           = @{ "A2" = "Waves_1663510155264" }
    "CheckPoint"            = @{ "A2" = "CheckPoint_1663510155424"; "A3" = "CheckPoint_1663510155521" }
    "Adobe"                 = @{ "A2" = "Adobe_1663510155457" }
    "Business Objects"      = @{ "A2" = "Business_1663510155555" }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
